$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 34 (pushes existing rows 34-54 down to 36-56),
# inheriting formatting (incl. the date style on column D) from the row above.
$ws.Range("A34:A35").EntireRow.Insert()

# New weekly row: Primera quality, week of 2021-11-29 (serial 44529)
$ws.Range("A34").Value = 9
$ws.Range("B34").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C34").Value = "Metropolitana"
$ws.Range("D34").Value = 44529
$ws.Range("E34").Value = 13
$ws.Range("F34").Value = 100114002
$ws.Range("G34").Value = "Camote"
$ws.Range("H34").Value = "Sin especificar"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 1010
$ws.Range("K34").Value = 12000
$ws.Range("L34").Value = 13000
$ws.Range("M34").Value = 12500
$ws.Range("N34").Value = '$/malla 18 kilos'
$ws.Range("O34").Value = "Perú"
$ws.Range("P34").Value = 694
$ws.Range("Q34").Value = 18
$ws.Range("R34").Value = "Hortaliza"

# New weekly row: Segunda quality, same week (serial 44529)
$ws.Range("A35").Value = 9
$ws.Range("B35").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C35").Value = "Metropolitana"
$ws.Range("D35").Value = 44529
$ws.Range("E35").Value = 13
$ws.Range("F35").Value = 100114002
$ws.Range("G35").Value = "Camote"
$ws.Range("H35").Value = "Sin especificar"
$ws.Range("I35").Value = "Segunda"
$ws.Range("J35").Value = 970
$ws.Range("K35").Value = 9000
$ws.Range("L35").Value = 10000
$ws.Range("M35").Value = 9495
$ws.Range("N35").Value = '$/malla 18 kilos'
$ws.Range("O35").Value = "Perú"
$ws.Range("P35").Value = 528
$ws.Range("Q35").Value = 18
$ws.Range("R35").Value = "Hortaliza"
